$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 59; this shifts rows 59..156 down to 60..157
# (Excel copies formatting from the row above, which keeps D's date style).
$ws.Rows.Item(59).Insert()

# Populate the new row 59 with the new data point.
$ws.Cells.Item(59, 1).Value = 8
$ws.Cells.Item(59, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = 44792
$ws.Cells.Item(59, 5).Value = 4
$ws.Cells.Item(59, 6).Value = 100112001
$ws.Cells.Item(59, 7).Value = "Berenjena"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 500
$ws.Cells.Item(59, 11).Value = 10000
$ws.Cells.Item(59, 12).Value = 11000
$ws.Cells.Item(59, 13).Value = 10500
$ws.Cells.Item(59, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(59, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(59, 16).Value = 210
$ws.Cells.Item(59, 17).Value = 50
$ws.Cells.Item(59, 18).Value = "Hortaliza"
